$wb = $excel.ActiveWorkbook

# --- Description sheet: text fixes -------------------------------------
$wsDesc = $wb.Worksheets.Item("Description")

# Add a missing full stop to the end of the "Desc Body" paragraph (B5)
$wsDesc.Range("B5").Value = "In 2014-15, 14.8 per cent of Australian adults were daily smokers. Nationally, between 2007-08 and 2014-15 there was a significant fall in the rate of smoking (4.3 percentage points), and significant falls in all states and territories except for the Northern Territory."

# Fix the mis-placed soft hyphens in the Australian Health Survey reference (B10)
$wsDesc.Range("B10").Value = "ABS (unpublished) Australian Health Survey 2011­-13 (2011-­12 core component)"

# The B5 paragraph now wraps slightly differently once the text changed, so
# nudge the row height back to what it auto-fit to.
$wsDesc.Rows.Item(5).RowHeight = 55.2

# --- Selections, matching where the user finished editing ---------------
$wsData = $wb.Worksheets.Item("Data")
$wsData.Activate()
$wsData.Range("A1").Select()

$wsDesc.Activate()
$wsDesc.Range("B10").Select()
